$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on Price (D) and Volume (E) columns for rows 2-51
# so numeric-looking strings (e.g. "1.007") are stored as text, matching the
# original inline-string cell type instead of being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '28.651.49'
$ws.Range("E2").Value = '  +2.27%  '

# Row 3
$ws.Range("D3").Value = '1.873.12'
$ws.Range("E3").Value = '  +2.44%  '

# Row 4
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  +0.41%  '

# Row 5
$ws.Range("D5").Value = '326.57'
$ws.Range("E5").Value = '  -0.41%  '

# Row 6
$ws.Range("E6").Value = '  +0.42%  '

# Row 7
$ws.Range("D7").Value = '0.4668'
$ws.Range("E7").Value = '  +1.05%  '

# Row 8
$ws.Range("D8").Value = '0.3893'
$ws.Range("E8").Value = '  +0.85%  '

# Row 9
$ws.Range("D9").Value = '0.07890'
$ws.Range("E9").Value = '  +0.30%  '

# Row 10
$ws.Range("D10").Value = '0.9738'
$ws.Range("E10").Value = '  +1.68%  '

# Row 11
$ws.Range("D11").Value = '22.02'
$ws.Range("E11").Value = '  +0.86%  '

# Row 12
$ws.Range("D12").Value = '1.918.78'
$ws.Range("E12").Value = '  +4.33%  '

# Row 13
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Value = '6.993'
$ws.Range("E13").Value = '  +1.68%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '5.711'
$ws.Range("E14").Value = '  +1.11%  '

# Row 15
$ws.Range("D15").Value = '0.06983'
$ws.Range("E15").Value = '  +3.19%  '

# Row 16
$ws.Range("D16").Value = '88.21'
$ws.Range("E16").Value = '  +1.86%  '

# Row 17
$ws.Range("E17").Value = '  +0.43%  '

# Row 18
$ws.Range("D18").Value = '0.00001005'
$ws.Range("E18").Value = '  +1.11%  '

# Row 19
$ws.Range("D19").Value = '16.85'
$ws.Range("E19").Value = '  +1.47%  '

# Row 20
$ws.Range("E20").Value = '  +0.34%  '

# Row 21
$ws.Range("D21").Value = '28.662.83'
$ws.Range("E21").Value = '  +2.24%  '

# Row 22
$ws.Range("D22").Value = '5.301'
$ws.Range("E22").Value = '  -0.01%  '

# Row 23
$ws.Range("D23").Value = '11.02'
$ws.Range("E23").Value = '  +0.42%  '

# Row 24
$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").Value = '2.165.64'
$ws.Range("E24").Value = '  +2.52%  '

# Row 25
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '2.116'
$ws.Range("E25").Value = '  +0.74%  '

# Row 26
$ws.Range("D26").Value = '152.53'
$ws.Range("E26").Value = '  -0.81%  '

# Row 27
$ws.Range("D27").Value = '19.24'
$ws.Range("E27").Value = '  +0.27%  '

# Row 28
$ws.Range("D28").Value = '5.760'
$ws.Range("E28").Value = '  +0.85%  '

# Row 29
$ws.Range("D29").Value = '1.989'
$ws.Range("E29").Value = '  +1.01%  '

# Row 30
$ws.Range("D30").Value = '119.31'
$ws.Range("E30").Value = '  +1.91%  '

# Row 31
$ws.Range("D31").Value = '0.09355'
$ws.Range("E31").Value = '  +1.13%  '

# Row 32
$ws.Range("D32").Value = '0.9201'
$ws.Range("E32").Value = '  -1.58%  '

# Row 33
$ws.Range("D33").Value = '5.275'
$ws.Range("E33").Value = '  -0.32%  '

# Row 34
$ws.Range("D34").Value = '1.340'
$ws.Range("E34").Value = '  +2.01%  '

# Row 35
$ws.Range("D35").Value = '3.347'
$ws.Range("E35").Value = '  +0.79%  '

# Row 36
$ws.Range("D36").Value = '0.05803'
$ws.Range("E36").Value = '  -1.04%  '

# Row 37
$ws.Range("D37").Value = '0.02097'
$ws.Range("E37").Value = '  -2.00%  '

# Row 38
$ws.Range("D38").Value = '1.143'
$ws.Range("E38").Value = '  +0.06%  '

# Row 39
$ws.Range("D39").Value = '7.779'
$ws.Range("E39").Value = '  +0.80%  '

# Row 40
$ws.Range("D40").Value = '0.5625'
$ws.Range("E40").Value = '  +0.87%  '

# Row 41
$ws.Range("D41").Value = '0.1785'
$ws.Range("E41").Value = '  +1.30%  '

# Row 42
$ws.Range("D42").Value = '9.797'
$ws.Range("E42").Value = '  -0.93%  '

# Row 43
$ws.Range("D43").Value = '0.07220'
$ws.Range("E43").Value = '  +3.03%  '

# Row 44
$ws.Range("D44").Value = '11.70'
$ws.Range("E44").Value = '  +1.36%  '

# Row 45
$ws.Range("D45").Value = '0.5317'
$ws.Range("E45").Value = '  +1.07%  '

# Row 46
$ws.Range("D46").Value = '1.166'
$ws.Range("E46").Value = '  -3.99%  '

# Row 47
$ws.Range("D47").Value = '1.826'
$ws.Range("E47").Value = '  -0.30%  '

# Row 48
$ws.Range("D48").Value = '113.37'
$ws.Range("E48").Value = '  +0.67%  '

# Row 49
$ws.Range("D49").Value = '2.062'
$ws.Range("E49").Value = '  -3.63%  '

# Row 50
$ws.Range("D50").Value = '2.364'
$ws.Range("E50").Value = '  +1.90%  '

# Row 51
$ws.Range("E51").Value = '  +0.47%  '
